$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.67
$ws.Range("M2").Value = 1.05
$ws.Range("N2").Value = 11
$ws.Range("O2").Value = 1.29
$ws.Range("P2").Value = 3.5
$ws.Range("Q2").Value = 1.98
$ws.Range("R2").Value = 1.83
